# Atualizei dados bibi e add
# Apply updated faturamento (billing) figures to rows 2-6 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bibi Cell Mundi)
$ws.Range("V2").Value = 15746.6
$ws.Range("AG2").Value = 184673.56

# Row 3 (Bibi Cell Vieiralves)
$ws.Range("T3").Value = 4390
$ws.Range("V3").Value = 5081
$ws.Range("W3").Value = 60
$ws.Range("AG3").Value = 90623.00999999999

# Row 4 (Bibi Cell Manauara)
$ws.Range("V4").Value = 2713
$ws.Range("AG4").Value = 64410.9

# Row 5 (Bibi Cell Ponta Negra)
$ws.Range("V5").Value = 2702.02
$ws.Range("AG5").Value = 57513.79

# Row 6 (total)
$ws.Range("T6").Value = 15007.65
$ws.Range("V6").Value = 26242.62
$ws.Range("W6").Value = 60
$ws.Range("AG6").Value = 397221.26
